# Commit Reports Fixed 03/09/2021
#
# Updates the test-data workbook's report sheets: refreshes the sample
# Start/End Date values used by the date-range reports, refreshes the
# "ShowDateRange" search-string rows (adding two more sample rows), fixes
# one Advanced Search sample value, and leaves the UI selection on the
# "ShowDateRange" tab, matching how the workbook was left after the edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Show (sheet1) - no data changes; just re-visit/re-select the cell that
# was active when the workbook was last saved.
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("Show")
$wsShow.Activate()
$wsShow.Range("D9").Select()

# ---------------------------------------------------------------------
# ShowInNewPage (sheet3) - no data changes.
# ---------------------------------------------------------------------
$wsShowInNewPage = $wb.Worksheets.Item("ShowInNewPage")
$wsShowInNewPage.Activate()
$wsShowInNewPage.Range("D29").Select()

# ---------------------------------------------------------------------
# ExportReport (sheet5) - no data changes, selection unchanged.
# ---------------------------------------------------------------------
$wsExportReport = $wb.Worksheets.Item("ExportReport")
$wsExportReport.Activate()
$wsExportReport.Range("C20").Select()

# ---------------------------------------------------------------------
# DrillDown (sheet9) - no data changes, selection unchanged.
# ---------------------------------------------------------------------
$wsDrillDown = $wb.Worksheets.Item("DrillDown")
$wsDrillDown.Activate()
$wsDrillDown.Range("H4").Select()

# ---------------------------------------------------------------------
# Queries (sheet7) - no data changes, selection unchanged; just visit it
# so its own state is untouched before we move the "active" tab away.
# ---------------------------------------------------------------------
$wsQueries = $wb.Worksheets.Item("Queries")
$wsQueries.Activate()
$wsQueries.Range("E2").Select()

# ---------------------------------------------------------------------
# ShowInNewPageDateRange (sheet4) - refresh the sample Start/End dates.
# ---------------------------------------------------------------------
$wsShowInNewPageDateRange = $wb.Worksheets.Item("ShowInNewPageDateRange")
$wsShowInNewPageDateRange.Activate()
$wsShowInNewPageDateRange.Range("D2").Value = "19-10-2020 00:00:00"
$wsShowInNewPageDateRange.Range("E2").Value = "20-10-2020 00:00:00"
$wsShowInNewPageDateRange.Range("D6").Select()

# ---------------------------------------------------------------------
# ExportReportDateRange (sheet6) - refresh the sample Start/End dates.
# ---------------------------------------------------------------------
$wsExportReportDateRange = $wb.Worksheets.Item("ExportReportDateRange")
$wsExportReportDateRange.Activate()
$wsExportReportDateRange.Range("D2").Value = "01-04-2020 00:00:00"
$wsExportReportDateRange.Range("E2").Value = "02-06-2020 00:00:00"
$wsExportReportDateRange.Range("D2").Select()

# ---------------------------------------------------------------------
# AdvanceSearch (sheet8) - fix row 2's "Search String2" sample value.
# ---------------------------------------------------------------------
$wsAdvanceSearch = $wb.Worksheets.Item("AdvanceSearch")
$wsAdvanceSearch.Activate()
$wsAdvanceSearch.Range("J2").Value = "Chat"
$wsAdvanceSearch.Range("E15").Select()

# ---------------------------------------------------------------------
# ShowDateRange (sheet2) - refresh the sample Start/End dates for the
# existing rows, fix up the search-string column, and append two more
# sample rows (5 and 6) so the sheet now covers rows 1-6.
# ---------------------------------------------------------------------
$wsShowDateRange = $wb.Worksheets.Item("ShowDateRange")
$wsShowDateRange.Activate()

$wsShowDateRange.Range("D2").Value = "19-10-2020 00:00:00"
$wsShowDateRange.Range("E2").Value = "20-10-2020 00:00:00"
$wsShowDateRange.Range("F2").Value = "30"

$wsShowDateRange.Range("D3").Value = "19-10-2020 00:00:00"
$wsShowDateRange.Range("E3").Value = "20-10-2020 00:00:00"
$wsShowDateRange.Range("F3").Value = "VoiceSkill1"

$wsShowDateRange.Range("D4").Value = "19-10-2020 00:00:00"
$wsShowDateRange.Range("E4").Value = "20-10-2020 00:00:00"
$wsShowDateRange.Range("F4").Value = "ice"

# New row 5, modelled on row 4's layout/format.
$wsShowDateRange.Range("A4:F4").Copy()
$wsShowDateRange.Range("A5:F5").PasteSpecial(-4122)
$wsShowDateRange.Range("A5").Value = "Agent"
$wsShowDateRange.Range("B5").Value = "OCM Skill Historical Report"
$wsShowDateRange.Range("C5").Value = "Date Range"
$wsShowDateRange.Range("D5").Value = "19-10-2020 00:00:00"
$wsShowDateRange.Range("E5").Value = "20-10-2020 00:00:00"
$wsShowDateRange.Range("F5").Value = "Email"

# New row 6, modelled on row 5's layout/format.
$wsShowDateRange.Range("A5:F5").Copy()
$wsShowDateRange.Range("A6:F6").PasteSpecial(-4122)
$wsShowDateRange.Range("A6").Value = "Agent"
$wsShowDateRange.Range("B6").Value = "OCM Skill Historical Report"
$wsShowDateRange.Range("C6").Value = "Date Range"
$wsShowDateRange.Range("D6").Value = "19-10-2020 00:00:00"
$wsShowDateRange.Range("E6").Value = "20-10-2020 00:00:00"
$wsShowDateRange.Range("F6").Value = "Skill11"

$excel.CutCopyMode = $false

$wsShowDateRange.Range("B9").Select()
